$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.901.98"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.603.04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.482"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0613"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.245"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.85"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.827.02"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.605.39"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.510"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.877.08"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.09"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0723"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.87%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "189.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.17"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.34"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.33%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.93"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.129"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.70"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.52"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.95"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.20"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0470"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.08"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.02"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.40"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.110.35"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.801"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.86%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.496"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "95.56"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.739.35"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.747"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.05"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0114"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "53.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.46"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0511"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.410"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.33"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.74%  "

Write-Output "Applied all changes"
